$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 109, shifting existing rows 109-113 down to 110-114
$ws.Rows.Item(109).Insert()

# Populate the new row 109 with the new record
$ws.Cells.Item(109, 1).Value = 9
$ws.Cells.Item(109, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(109, 3).Value = "Metropolitana"
$ws.Cells.Item(109, 4).Value = 44931
$ws.Cells.Item(109, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(109, 5).Value = 13
$ws.Cells.Item(109, 6).Value = "Fruta"
$ws.Cells.Item(109, 7).Value = 100101
$ws.Cells.Item(109, 8).Value = "Berries"
$ws.Cells.Item(109, 9).Value = 100101004
$ws.Cells.Item(109, 10).Value = "Frambuesa"
$ws.Cells.Item(109, 11).Value = "Sin especificar"
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 400
$ws.Cells.Item(109, 14).Value = 8000
$ws.Cells.Item(109, 15).Value = 8000
$ws.Cells.Item(109, 16).Value = 8000
$ws.Cells.Item(109, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(109, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(109, 19).Value = 4000
$ws.Cells.Item(109, 20).Value = 2
